# Update "想去人数" (F column) figures across all sheets to reflect the
# regenerated scrape output (gh-pages build at 456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 12857
$ws.Range("F3").Value = 7189
$ws.Range("F6").Value = 454
$ws.Range("F11").Value = 148
$ws.Range("F15").Value = 69
$ws.Range("F16").Value = 1020
$ws.Range("F18").Value = 250
$ws.Range("F19").Value = 372
$ws.Range("F21").Value = 282
$ws.Range("F23").Value = 50
$ws.Range("F24").Value = 181
$ws.Range("F25").Value = 376
$ws.Range("F26").Value = 5247
$ws.Range("F27").Value = 72
$ws.Range("F28").Value = 1435
$ws.Range("F30").Value = 1396
$ws.Range("F31").Value = 67
$ws.Range("F32").Value = 54
$ws.Range("F33").Value = 1370
$ws.Range("F35").Value = 6
$ws.Range("F36").Value = 599
$ws.Range("F38").Value = 3741

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 3744
$ws.Range("F5").Value = 3744
$ws.Range("F8").Value = 60

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 9288
$ws.Range("F3").Value = 563
$ws.Range("F4").Value = 2029

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 9288
$ws.Range("F3").Value = 563
$ws.Range("F4").Value = 2029
$ws.Range("F5").Value = 12857
$ws.Range("F6").Value = 7189
$ws.Range("F8").Value = 3744
$ws.Range("F9").Value = 454
$ws.Range("F11").Value = 148
$ws.Range("F15").Value = 69
$ws.Range("F16").Value = 1020
$ws.Range("F18").Value = 250
$ws.Range("F19").Value = 372
$ws.Range("F21").Value = 282
$ws.Range("F23").Value = 50
$ws.Range("F27").Value = 181
$ws.Range("F28").Value = 376
$ws.Range("F29").Value = 5247
$ws.Range("F30").Value = 72
$ws.Range("F31").Value = 1435
$ws.Range("F36").Value = 1396
$ws.Range("F37").Value = 67
$ws.Range("F38").Value = 1370
$ws.Range("F40").Value = 599
$ws.Range("F47").Value = 3741
